$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18 ("Passes"): fix wording in the "Note:" paragraph of the body text
#   "compiler pass" -> "compiler"
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$body18 = $s18.Shapes.Item(4)
$tr18 = $body18.TextFrame.TextRange
$para18 = $tr18.Paragraphs(3)
# Clear to unrelated text first so the new text is written as a single run
# instead of being diffed/split against the old wording.
$para18.Text = "x"
$tr18.Paragraphs(3).Text = "Note: Some authors restrict the definition of compiler to a traversal that involves disk I/O, but we will use a more general definition."

# ---------------------------------------------------------------------------
# Slide 19 ("Single-pass Versus Multi-pass Compilers"): reword bullet
#   "can exploit concurrency and multiprocessor architectures"
#   -> "ideal for multiprocessor systems"
# ---------------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$body19 = $s19.Shapes.Item(4)
$tr19 = $body19.TextFrame.TextRange
$para19 = $tr19.Paragraphs(5)
$para19.Text = "x"
$tr19.Paragraphs(5).Text = "ideal for multiprocessor systems"

# ---------------------------------------------------------------------------
# Slide 7 ("Parser (Syntax Analysis)"): nudge two shapes that are part of the
# parse-tree diagram so the "intLiteral1" label and its connector line up
# with the widened text box.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# "Text Box 19" (intLiteral1 label)
$intLiteralBox = $s7.Shapes.Item(12)
$intLiteralBox.Left = 421.5911811023622
$intLiteralBox.Width = 84.56795275590551

# "AutoShape 23" (connector leading into the label)
$connector = $s7.Shapes.Item(16)
$connector.Width = 47.7502

# ---------------------------------------------------------------------------
# Slide 9 ("Code Generator"): clarify when the code generator is referred to
# as "intermediate"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$body9 = $s9.Shapes.Item(4)
$tr9 = $body9.TextFrame.TextRange
$para9 = $tr9.Paragraphs(3)
$para9.Text = "x"
$tr9.Paragraphs(3).Text = "If the low-level representation is assembly language or if it is machine independent, then this component of the compiler is often referred to as an “intermediate” code generator."
